# Updated cryptos list on Sat Sep  7 21:53:19 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# crypto table, and fixes the Uniswap / BitcoinCash rows which had swapped
# places (row 20 is now BitcoinCash, row 21 is now Uniswap) together with
# their refreshed price/volume figures.
#
# Several "Price" values look numeric (single decimal point, e.g. "492.43")
# but must stay literal text (matching the original inlineStr cells), so we
# force the cell to Text format before writing those values; otherwise Excel
# would silently convert them to real numbers (and drop things like the
# trailing zero in "0.370" / "6.40" / "17.60").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "53.938.96"
$ws.Range("E2").Value = "  +1.03%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.247.16"
$ws.Range("E3").Value = "  +2.66%  "

# Row 4 - TetherUSD (price unchanged)
$ws.Range("E4").Value = "  +0.13%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "492.43"
$ws.Range("E5").Value = "  +2.21%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "127.19"
$ws.Range("E6").Value = "  +2.64%  "

# Row 7 - USDC (price unchanged)
$ws.Range("E7").Value = "  +0.04%  "

# Row 8 - XRP
Set-TextValue $ws.Range("D8") "0.527"
$ws.Range("E8").Value = "  +1.72%  "

# Row 9 - Dogecoin (price unchanged)
$ws.Range("E9").Value = "  +4.31%  "

# Row 10 - TRON (price unchanged)
$ws.Range("E10").Value = "  +2.62%  "

# Row 11 - Cardano (price unchanged)
$ws.Range("E11").Value = "  +4.00%  "

# Row 12 - Toncoin (price unchanged)
$ws.Range("E12").Value = "  +1.74%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D13") "2.648.89"
$ws.Range("E13").Value = "  +2.97%  "

# Row 14 - Avalanche (price unchanged)
$ws.Range("E14").Value = "  +3.79%  "

# Row 15 - WrappedBTC
Set-TextValue $ws.Range("D15") "53.882.33"
$ws.Range("E15").Value = "  +1.11%  "

# Row 16 - ShibaInu
Set-TextValue $ws.Range("D16") "0.0000129"
$ws.Range("E16").Value = "  +1.36%  "

# Row 17 - WrappedEther
Set-TextValue $ws.Range("D17") "2.255.53"
$ws.Range("E17").Value = "  +2.35%  "

# Row 18 - Chainlink (price unchanged)
$ws.Range("E18").Value = "  +5.03%  "

# Row 20 - now BitcoinCash (was Uniswap)
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D20") "299.29"
$ws.Range("E20").Value = "  +2.22%  "

# Row 21 - now Uniswap (was BitcoinCash)
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D21") "6.40"
$ws.Range("E21").Value = "  +5.99%  "

# Row 22 - Dai (price unchanged)
$ws.Range("E22").Value = "  -0.08%  "

# Row 23 - LEO (price unchanged)
$ws.Range("E23").Value = "  -2.28%  "

# Row 24 - Litecoin (price unchanged)
$ws.Range("E24").Value = "  -1.03%  "

# Row 25 - Binance-PegBSC-USD (price unchanged)
$ws.Range("E25").Value = "  +2.16%  "

# Row 26 - Polygon
Set-TextValue $ws.Range("D26") "0.370"
$ws.Range("E26").Value = "  +1.85%  "

# Row 27 - WrappedeETH
Set-TextValue $ws.Range("D27") "2.355.66"
$ws.Range("E27").Value = "  +2.99%  "

# Row 28 - Kaspa (price unchanged)
$ws.Range("E28").Value = "  +3.26%  "

# Row 29 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D29") "7.02"
$ws.Range("E29").Value = "  +1.08%  "

# Row 30 - Monero
Set-TextValue $ws.Range("D30") "167.31"
$ws.Range("E30").Value = "  +1.51%  "

# Row 31 - PancakeSwap (price unchanged)
$ws.Range("E31").Value = "  +2.32%  "

# Row 32 - PEPE (price unchanged)
$ws.Range("E32").Value = "  +4.01%  "

# Row 33 - Aptos
Set-TextValue $ws.Range("D33") "5.83"
$ws.Range("E33").Value = "  +3.20%  "

# Row 34 - USDe (price unchanged)
$ws.Range("E34").Value = "  +0.12%  "

# Row 35 - FirstDigitalUSD (price unchanged)
$ws.Range("E35").Value = "  +0.30%  "

# Row 36 - Fetch.AI (price unchanged)
$ws.Range("E36").Value = "  +1.42%  "

# Row 37 - EthereumClassic
Set-TextValue $ws.Range("D37") "17.60"
$ws.Range("E37").Value = "  +2.29%  "

# Row 38 - SuiNetwork
Set-TextValue $ws.Range("D38") "0.899"
$ws.Range("E38").Value = "  +10.35%  "

# Row 39 - ImmutableX (price unchanged)
$ws.Range("E39").Value = "  +3.34%  "

# Row 40 - NEARProtocol
Set-TextValue $ws.Range("D40") "3.66"
$ws.Range("E40").Value = "  +4.25%  "

# Row 41 - OKB
Set-TextValue $ws.Range("D41") "35.67"
$ws.Range("E41").Value = "  +0.00%  "

# Row 42 - Stacks (price unchanged)
$ws.Range("E42").Value = "  +3.38%  "

# Row 43 - PolygonEcosystemToken (price unchanged)
$ws.Range("E43").Value = "  +2.16%  "

# Row 44 - Filecoin (price unchanged)
$ws.Range("E44").Value = "  +3.88%  "

# Row 45 - RenderToken
Set-TextValue $ws.Range("D45") "4.91"
$ws.Range("E45").Value = "  +4.49%  "

# Row 46 - Aave
Set-TextValue $ws.Range("D46") "124.98"
$ws.Range("E46").Value = "  +1.87%  "

# Row 47 - Stellar (price unchanged)
$ws.Range("E47").Value = "  +1.64%  "

# Row 48 - Mantle
Set-TextValue $ws.Range("D48") "0.540"
$ws.Range("E48").Value = "  +1.97%  "

# Row 49 - Bittensor
Set-TextValue $ws.Range("D49") "236.16"
$ws.Range("E49").Value = "  +3.80%  "

# Row 50 - Hedera (price unchanged)
$ws.Range("E50").Value = "  +3.44%  "

# Row 51 - VeChain
Set-TextValue $ws.Range("D51") "0.0201"
$ws.Range("E51").Value = "  +1.93%  "
